# Update the dSF column (F) values on the active worksheet to reflect the
# repulled/recalculated data (repull data, push all data, mean calculation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 6
    6  = -3
    7  = -6
    8  = -2
    10 = -4
    11 = 2
    12 = 4
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = -2
    19 = -5
    21 = -2
    22 = -3
    23 = -1
    24 = 1
    25 = 1
    27 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
